$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A92").Value = "GRT-USD"
